$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column D (already styled, empty) gets the "Rien pour le moment..." value
# and a new column E cell with "ignore" is added.
$rowsWithExistingStyle = @(11, 13, 21, 23, 45, 46)
foreach ($r in $rowsWithExistingStyle) {
    $ws.Range("D$r").Value = "Rien pour le moment…"
    $ws.Range("E$r").Value = "ignore"
}

# Row 58: column D has no cell at all yet, so copy C58's format (style) first,
# then set the value; then add column E.
$ws.Range("C58").Copy()
$ws.Range("D58").PasteSpecial(-4122)
$ws.Range("D58").Value = "Rien pour le moment…"
$ws.Range("E58").Value = "ignore"

# Update the active selection to reflect the saved view state (F57 -> F58).
$ws.Activate()
$ws.Range("F58").Select()
